$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlinks on column A (members used to be linked by email)
$ws.Hyperlinks.Delete()

# Header: email do membro -> Id do membro
$ws.Range("A1").Value = "Id do membro"

# Replace the member e-mail addresses in column A with numeric member ids
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 3
$ws.Range("A7").Value = 3

# Re-apply the bordered/centered look (matches the style already used by column B)
# by copying its format onto the new id column, instead of the old hyperlink style.
$fmtSrc = $ws.Range("B2")
$fmtSrc.Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Selection moved in the saved file
$ws.Range("A12").Select()
